$d = $word.ActiveDocument

# Phase 1: replace each original expression with a unique, collision-free placeholder
# token so that no intermediate result can accidentally be re-matched by a later rule
# (e.g. 14÷3= becomes 31÷4=, which is itself also one of the values being replaced).
$d.Content.Find.Execute("26÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER0§", 2)
$d.Content.Find.Execute("43÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER1§", 2)
$d.Content.Find.Execute("36÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER2§", 2)
$d.Content.Find.Execute("61÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER3§", 2)
$d.Content.Find.Execute("22÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER4§", 2)
$d.Content.Find.Execute("62÷2=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER5§", 2)
$d.Content.Find.Execute("31÷4=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER6§", 2)
$d.Content.Find.Execute("84÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER7§", 2)
$d.Content.Find.Execute("61÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER8§", 2)
$d.Content.Find.Execute("69÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER9§", 2)
$d.Content.Find.Execute("38÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER10§", 2)
$d.Content.Find.Execute("66÷7=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER11§", 2)
$d.Content.Find.Execute("99÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER12§", 2)
$d.Content.Find.Execute("35÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER13§", 2)
$d.Content.Find.Execute("14÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER14§", 2)
$d.Content.Find.Execute("34÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER15§", 2)
$d.Content.Find.Execute("39÷8=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER16§", 2)
$d.Content.Find.Execute("12÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER17§", 2)
$d.Content.Find.Execute("60÷5=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER18§", 2)
$d.Content.Find.Execute("67÷6=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER19§", 2)
$d.Content.Find.Execute("15÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER20§", 2)
$d.Content.Find.Execute("23÷3=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER21§", 2)
$d.Content.Find.Execute("55÷2=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER22§", 2)
$d.Content.Find.Execute("91÷2=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER23§", 2)
$d.Content.Find.Execute("20÷9=", $false, $false, $false, $false, $false, $true, 1, $false, "§PLACEHOLDER24§", 2)

# Phase 2: replace the placeholders with the final new expressions
$d.Content.Find.Execute("§PLACEHOLDER0§", $false, $false, $false, $false, $false, $true, 1, $false, "98÷7=", 2)
$d.Content.Find.Execute("§PLACEHOLDER1§", $false, $false, $false, $false, $false, $true, 1, $false, "97÷5=", 2)
$d.Content.Find.Execute("§PLACEHOLDER2§", $false, $false, $false, $false, $false, $true, 1, $false, "33÷2=", 2)
$d.Content.Find.Execute("§PLACEHOLDER3§", $false, $false, $false, $false, $false, $true, 1, $false, "75÷5=", 2)
$d.Content.Find.Execute("§PLACEHOLDER4§", $false, $false, $false, $false, $false, $true, 1, $false, "94÷7=", 2)
$d.Content.Find.Execute("§PLACEHOLDER5§", $false, $false, $false, $false, $false, $true, 1, $false, "15÷7=", 2)
$d.Content.Find.Execute("§PLACEHOLDER6§", $false, $false, $false, $false, $false, $true, 1, $false, "81÷6=", 2)
$d.Content.Find.Execute("§PLACEHOLDER7§", $false, $false, $false, $false, $false, $true, 1, $false, "87÷8=", 2)
$d.Content.Find.Execute("§PLACEHOLDER8§", $false, $false, $false, $false, $false, $true, 1, $false, "99÷7=", 2)
$d.Content.Find.Execute("§PLACEHOLDER9§", $false, $false, $false, $false, $false, $true, 1, $false, "16÷3=", 2)
$d.Content.Find.Execute("§PLACEHOLDER10§", $false, $false, $false, $false, $false, $true, 1, $false, "70÷7=", 2)
$d.Content.Find.Execute("§PLACEHOLDER11§", $false, $false, $false, $false, $false, $true, 1, $false, "85÷9=", 2)
$d.Content.Find.Execute("§PLACEHOLDER12§", $false, $false, $false, $false, $false, $true, 1, $false, "64÷4=", 2)
$d.Content.Find.Execute("§PLACEHOLDER13§", $false, $false, $false, $false, $false, $true, 1, $false, "39÷9=", 2)
$d.Content.Find.Execute("§PLACEHOLDER14§", $false, $false, $false, $false, $false, $true, 1, $false, "31÷4=", 2)
$d.Content.Find.Execute("§PLACEHOLDER15§", $false, $false, $false, $false, $false, $true, 1, $false, "52÷3=", 2)
$d.Content.Find.Execute("§PLACEHOLDER16§", $false, $false, $false, $false, $false, $true, 1, $false, "57÷6=", 2)
$d.Content.Find.Execute("§PLACEHOLDER17§", $false, $false, $false, $false, $false, $true, 1, $false, "94÷2=", 2)
$d.Content.Find.Execute("§PLACEHOLDER18§", $false, $false, $false, $false, $false, $true, 1, $false, "65÷3=", 2)
$d.Content.Find.Execute("§PLACEHOLDER19§", $false, $false, $false, $false, $false, $true, 1, $false, "51÷4=", 2)
$d.Content.Find.Execute("§PLACEHOLDER20§", $false, $false, $false, $false, $false, $true, 1, $false, "31÷7=", 2)
$d.Content.Find.Execute("§PLACEHOLDER21§", $false, $false, $false, $false, $false, $true, 1, $false, "82÷8=", 2)
$d.Content.Find.Execute("§PLACEHOLDER22§", $false, $false, $false, $false, $false, $true, 1, $false, "18÷4=", 2)
$d.Content.Find.Execute("§PLACEHOLDER23§", $false, $false, $false, $false, $false, $true, 1, $false, "79÷3=", 2)
$d.Content.Find.Execute("§PLACEHOLDER24§", $false, $false, $false, $false, $false, $true, 1, $false, "74÷3=", 2)
